$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = 'sd'
$ws.Range("J3").Value = 'Statement-non-opinion'
$ws.Range("I4").Value = 'sd'
$ws.Range("J4").Value = 'Statement-non-opinion'
$ws.Range("I5").Value = 'sd'
$ws.Range("J5").Value = 'Statement-non-opinion'
$ws.Range("I8").Value = 'sv'
$ws.Range("J8").Value = 'Statement-opinion'
$ws.Range("I10").Value = 'sv'
$ws.Range("J10").Value = 'Statement-opinion'
$ws.Range("I11").Value = 'sd'
$ws.Range("J11").Value = 'Statement-non-opinion'
$ws.Range("I14").Value = 'sd'
$ws.Range("J14").Value = 'Statement-non-opinion'
$ws.Range("I19").Value = 'sd'
$ws.Range("J19").Value = 'Statement-non-opinion'
$ws.Range("I44").Value = 'aa'
$ws.Range("J44").Value = 'Agree/Accept'
$ws.Range("I50").Value = 'sd'
$ws.Range("J50").Value = 'Statement-non-opinion'
$ws.Range("I51").Value = 'sd'
$ws.Range("J51").Value = 'Statement-non-opinion'
$ws.Range("I59").Value = 'sv'
$ws.Range("J59").Value = 'Statement-opinion'
$ws.Range("I60").Value = 'ba'
$ws.Range("J60").Value = 'Appreciation'
$ws.Range("I71").Value = 'sd'
$ws.Range("J71").Value = 'Statement-non-opinion'
$ws.Range("I73").Value = 'sd'
$ws.Range("J73").Value = 'Statement-non-opinion'
$ws.Range("I74").Value = 'sd'
$ws.Range("J74").Value = 'Statement-non-opinion'
$ws.Range("I75").Value = 'sv'
$ws.Range("J75").Value = 'Statement-opinion'
$ws.Range("I84").Value = 'sd'
$ws.Range("J84").Value = 'Statement-non-opinion'
$ws.Range("I97").Value = 'sv'
$ws.Range("J97").Value = 'Statement-opinion'
$ws.Range("I98").Value = '%'
$ws.Range("J98").Value = 'Uninterpretable'
$ws.Range("I109").Value = 'sd'
$ws.Range("J109").Value = 'Statement-non-opinion'
$ws.Range("I110").Value = 'b'
$ws.Range("J110").Value = 'Acknowledge (Backchannel)'
$ws.Range("I114").Value = 'ba'
$ws.Range("J114").Value = 'Appreciation'
$ws.Range("I116").Value = 'b'
$ws.Range("J116").Value = 'Acknowledge (Backchannel)'
$ws.Range("I132").Value = 'sv'
$ws.Range("J132").Value = 'Statement-opinion'
